$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: for values that look like plain numbers (e.g. "536.94", "70.00", "1.20"),
# Excels Range.Value setter auto-converts the text to a real number and
# trailing zeros / precision get lost. Forcing the cell to Text format first
# (then resetting the style back to Normal so no visible formatting changes)
# keeps these as literal text, matching the source data.

$ws.Range('D2').Value = '60.009.58'
$ws.Range('E2').Value = '  +2.51%  '
$ws.Range('D3').Value = '3.203.82'
$ws.Range('E3').Value = '  +1.40%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '536.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.10'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.61%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.530'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.81%  '
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('E10').Value = '  +2.81%  '
$ws.Range('E11').Value = '  +2.68%  '
$ws.Range('D12').Value = '3.754.19'
$ws.Range('E12').Value = '  +1.37%  '
$ws.Range('E13').Value = '  -1.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.89'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('E15').Value = '  +1.48%  '
$ws.Range('D16').Value = '60.074.99'
$ws.Range('E16').Value = '  +2.53%  '
$ws.Range('D17').Value = '3.202.80'
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.27'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.17'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('E20').Value = '  +0.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '375.47'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.73%  '
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('E23').Value = '  +1.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('E25').Value = '  +1.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.75'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +7.72%  '
$ws.Range('E27').Value = '  +0.31%  '
$ws.Range('D28').Value = '0.0₃0897'
$ws.Range('E28').Value = '  +2.81%  '
$ws.Range('E29').Value = '  +0.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.34'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.15'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.38'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.46%  '
$ws.Range('B33').Value = 'Aptos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.67'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.93%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.20'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '156.83'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.26%  '
$ws.Range('E36').Value = '  -1.99%  '
$ws.Range('D37').Value = '2.801.45'
$ws.Range('E37').Value = '  +6.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.66'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0705'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('E41').Value = '  +1.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.88'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.96%  '
$ws.Range('E43').Value = '  +3.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.718'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.106'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.81%  '
$ws.Range('D46').Value = '3.246.00'
$ws.Range('E46').Value = '  +1.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.814'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +7.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.983'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('E49').Value = '  -0.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.61'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.68%  '
$ws.Range('E51').Value = '  +0.01%  '
